$wb = $excel.ActiveWorkbook

$wsK = $wb.Worksheets.Item("k_constants_log10")
$wsK.Range("A2").Value = 2
$wsK.Range("A3").Value = 2
$wsK.Range("A4").Select()

$wsI = $wb.Worksheets.Item("individual_shifts")
$wsI.Range("B2").Value = 8.4596999999999998
$wsI.Range("B3").Select()

$wsC = $wb.Worksheets.Item("chemical_shifts")
$wsC.Range("D2").Value = 8.4596999999999998
$wsC.Activate()
$wsC.Range("D3").Select()
